$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.745259523391724
$ws.Range("B1").Value = 2.630134105682373
$ws.Range("C1").Value = 2.617287397384644
$ws.Range("D1").Value = 2.701791048049927
$ws.Range("E1").Value = 3.400125503540039
